# peakList accessors for MSPeakListsSet
# Adds "X" marks to the "mslists" sheet for the peakLists / averagedPeakLists
# rows in the newly-supported columns (F = ionize, G = done).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mslists")

# Mark the new accessor support cells with "X"
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "X"

$ws.Range("G10").Value = "X"

$ws.Range("G11").Value = "X"

$ws.Range("G12").Value = "X"

$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"

# Update the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("G15").Select()
